$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("A1").Value = 1.483322143554688
$ws.Range("B1").Value = 3.64454460144043
$ws.Range("C1").Value = 6.022371768951416
$ws.Range("D1").Value = 1.477355718612671
$ws.Range("E1").Value = 0.8647710084915161
